$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before C; old C:F data shifts right to E:H
$ws.Columns("C:D").Insert()

# Populate the newly inserted C:D columns (duplicates of A:B header/first row)
$ws.Range("C1").Value = "CompCode"
$ws.Range("D1").Value = "WrkGrp"
$ws.Range("C2").Value = "'01"
$ws.Range("D2").Value = "COMP"

# Update the data values that changed in the (now shifted) columns F, G, H
$ws.Range("F2").Value = 104019
$ws.Range("G2").Value = "OH"
$ws.Range("H2").Value = 2

# Give the new block C1:H2 a thin box border around every cell
$ws.Range("C1:H2").Borders.LineStyle = 1
$ws.Range("C1:H2").Borders.Weight = 2

# Best-effort autofit of the new CompCode column width (closest achievable value)
$ws.Columns("C:C").ColumnWidth = 9.833333333333334

# Restore the saved selection/active cell
$ws.Range("H14").Select() | Out-Null
